# Implements: "Implemented getting number of lines for methods and classes."
#
# 1. Swap the `id`/`log` field rows on the classFields sheet (OrderController),
#    mirroring the reordering of those fields in sharedStrings.
# 2. Add two new worksheets at the end of the workbook:
#       - classNumberOfLines  (Class Name | Number of Lines)
#       - methodNumberOfLines (Class Name | Method Signature | Number of Lines)

$wb = $excel.ActiveWorkbook

# --- 1. classFields: swap the `id` and `log` field entries for OrderController ---
$fields = $wb.Worksheets.Item("classFields")
$fields.Range("B3").Value = "log"
$fields.Range("D3").Value = "org.slf4j.Logger"
$fields.Range("B5").Value = "id"
$fields.Range("D5").Value = "java.util.concurrent.atomic.AtomicLong"

# --- helper: write a text value into a cell without it being coerced to a number ---
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    if ($text -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $text
}

# --- 2a. classNumberOfLines ---
$classLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$classSheet = $wb.Worksheets.Add($null, $classLast)
$classSheet.Name = "classNumberOfLines"

$classRows = @(
    @("Class Name", "Number of Lines"),
    @("com.zatribune.spring.ecommerce.orders.controller.OrderController", "28"),
    @("com.zatribune.spring.ecommerce.orders.service.OrderService", "4"),
    @("com.zatribune.spring.ecommerce.orders.OrderApplicationTests", "5"),
    @("com.zatribune.spring.ecommerce.orders.service.OrderServiceImpl", "18"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "47"),
    @("com.zatribune.spring.ecommerce.orders.OrderApplication", "6")
)

$r = 1
foreach ($row in $classRows) {
    Set-TextCell $classSheet $r 1 $row[0]
    Set-TextCell $classSheet $r 2 $row[1]
    $r += 1
}

# --- 2b. methodNumberOfLines ---
$methodLast = $wb.Worksheets.Item($wb.Worksheets.Count)
$methodSheet = $wb.Worksheets.Add($null, $methodLast)
$methodSheet.Name = "methodNumberOfLines"

$methodRows = @(
    @("Class Name", "Method Signature", "Number of Lines"),
    @("com.zatribune.spring.ecommerce.orders.controller.OrderController", "create(domain.Order)", "1"),
    @("com.zatribune.spring.ecommerce.orders.controller.OrderController", "all()", "7"),
    @("com.zatribune.spring.ecommerce.orders.controller.OrderController", 'lambda$all$0(java.util.List, org.apache.kafka.streams.KeyValue)', "7"),
    @("com.zatribune.spring.ecommerce.orders.service.OrderService", "confirm(domain.Order, domain.Order)", "1"),
    @("com.zatribune.spring.ecommerce.orders.OrderApplicationTests", "contextLoads()", "2"),
    @("com.zatribune.spring.ecommerce.orders.service.OrderServiceImpl", "confirm(domain.Order, domain.Order)", "5"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "orders()", "3"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "paymentTopic()", "3"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "stockTopic()", "3"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "stream(org.apache.kafka.streams.StreamsBuilder)", "10"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "table(org.apache.kafka.streams.StreamsBuilder)", "5"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", "taskExecutor()", "8"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", 'lambda$table$1(java.lang.Long, domain.Order)', "5"),
    @("com.zatribune.spring.ecommerce.orders.config.KafkaConfig", 'lambda$stream$0(java.lang.Long, domain.Order)', "10"),
    @("com.zatribune.spring.ecommerce.orders.OrderApplication", "main(java.lang.String[])", "3")
)

$r = 1
foreach ($row in $methodRows) {
    Set-TextCell $methodSheet $r 1 $row[0]
    Set-TextCell $methodSheet $r 2 $row[1]
    Set-TextCell $methodSheet $r 3 $row[2]
    $r += 1
}
